# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 7197
$wsExhibition.Range("F5").Value = 150
$wsExhibition.Range("F6").Value = 1102
$wsExhibition.Range("F7").Value = 179
$wsExhibition.Range("F9").Value = 70
$wsExhibition.Range("F10").Value = 14

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7197
$wsAll.Range("F5").Value = 150
$wsAll.Range("F6").Value = 1102
$wsAll.Range("F7").Value = 179
$wsAll.Range("F10").Value = 70
$wsAll.Range("F11").Value = 14
